# Apply "Szenario 0" re-weighting to the Gewichtung sheet and tidy up
# rendering artifacts (total_weight column, helper styles, active sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gewichtung")

# --- Update group_weight (C) and within_group_weight (F) values ---------

# umweltbelastung
$ws.Range("C2").Value = 0.45
$ws.Range("F2").Value = 0.5
$ws.Range("C3").Value = 0.45
$ws.Range("F3").Value = 0.5

# langlebigkeit_wirtschaftlichkeit
$ws.Range("C4").Value = 0.05
$ws.Range("F4").Value = 0.5
$ws.Range("C5").Value = 0.05
$ws.Range("F5").Value = 0.25
$ws.Range("C6").Value = 0.05
$ws.Range("F6").Value = 0.25

# multifunktionale_nutzungsqualitaet
$ws.Range("C7").Value = 0.05
$ws.Range("F7").Value = 0.25
$ws.Range("C8").Value = 0.05
$ws.Range("F8").Value = 0.25
$ws.Range("C9").Value = 0.05
$ws.Range("F9").Value = 0.25
$ws.Range("C10").Value = 0.05
$ws.Range("F10").Value = 0.25

# kreislauffaehigkeit
$ws.Range("C11").Value = 0.45
$ws.Range("F11").Formula = "=1/3"
$ws.Range("C12").Value = 0.45
$ws.Range("F12").Formula = "=1/3"
$ws.Range("C13").Value = 0.45
$ws.Range("F13").Formula = "=1/3"

# --- Drop the total_weight column (G) ------------------------------------
# Clear the header text (keep its border/alignment style) and remove all
# the computed values below it entirely.
$ws.Range("G1:G13").ClearContents()

# --- Remove the now-unused column F helper style (fill-highlight) -------
# The F column used to carry a distinct "highlighted" style (col-level and
# per-cell); strip it back to the plain default and restore the header
# cell's normal bold/border/centered look.
$ws.Columns.Item(6).ClearFormats()
$f1 = $ws.Range("F1")
$f1.Borders.LineStyle = 1
$f1.Font.Bold = $true
$f1.VerticalAlignment = -4160
$f1.HorizontalAlignment = -4108

# --- Sheet selection / active sheet --------------------------------------
# "Gewichtung" becomes the active/visible tab, with a fresh selection.
$ws.Activate()
$ws.Range("K14").Select()
